# Natmi following Dr Hou advice
# Update Fn1-Mag LR-pairs sheet: recompute stats using min cell count of 3 (was 1),
# and expand target clusters from {M2} to {M2, sCs} for every sending cluster
# (ECs, FAPs, M2, sCs), producing 8 data rows (rows 2-9) instead of 4 (rows 2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Mag"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3133113333333333
$ws.Range("N2").Value = 0.939934
$ws.Range("O2").Value = 0.4010297802586483
$ws.Range("P2").Value = 0.4010297802586483
$ws.Range("Q2").Value = 12.76849937650178
$ws.Range("R2").Value = 114.916494388516
$ws.Range("S2").Value = 0.008528045179168314
$ws.Range("T2").Value = 0.008528045179168314

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Mag"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4679556666666667
$ws.Range("N3").Value = 1.403867
$ws.Range("O3").Value = 0.5989702197413518
$ws.Range("P3").Value = 0.5989702197413518
$ws.Range("Q3").Value = 19.07078041031756
$ws.Range("R3").Value = 171.637023692858
$ws.Range("S3").Value = 0.01273732113270025
$ws.Range("T3").Value = 0.01273732113270025

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Mag"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1689.289306666667
$ws.Range("H4").Value = 5067.86792
$ws.Range("I4").Value = 0.8814813868902838
$ws.Range("J4").Value = 0.8814813868902838
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3133113333333333
$ws.Range("N4").Value = 0.939934
$ws.Range("O4").Value = 0.4010297802586483
$ws.Range("P4").Value = 0.4010297802586483
$ws.Range("Q4").Value = 529.2734850574756
$ws.Range("R4").Value = 4763.46136551728
$ws.Range("S4").Value = 0.353500286886699
$ws.Range("T4").Value = 0.353500286886699

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Mag"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1689.289306666667
$ws.Range("H5").Value = 5067.86792
$ws.Range("I5").Value = 0.8814813868902838
$ws.Range("J5").Value = 0.8814813868902838
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4679556666666667
$ws.Range("N5").Value = 1.403867
$ws.Range("O5").Value = 0.5989702197413518
$ws.Range("P5").Value = 0.5989702197413518
$ws.Range("Q5").Value = 790.512503694071
$ws.Range("R5").Value = 7114.612533246639
$ws.Range("S5").Value = 0.5279811000035848
$ws.Range("T5").Value = 0.5279811000035848

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Mag"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 100.9654023333333
$ws.Range("H6").Value = 302.896207
$ws.Range("I6").Value = 0.05268435816499466
$ws.Range("J6").Value = 0.05268435816499466
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3133113333333333
$ws.Range("N6").Value = 0.939934
$ws.Range("O6").Value = 0.4010297802586483
$ws.Range("P6").Value = 0.4010297802586483
$ws.Range("Q6").Value = 31.63360482559311
$ws.Range("R6").Value = 284.702443430338
$ws.Range("S6").Value = 0.02112799657797573
$ws.Range("T6").Value = 0.02112799657797573

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Mag"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 100.9654023333333
$ws.Range("H7").Value = 302.896207
$ws.Range("I7").Value = 0.05268435816499466
$ws.Range("J7").Value = 0.05268435816499466
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4679556666666667
$ws.Range("N7").Value = 1.403867
$ws.Range("O7").Value = 0.5989702197413518
$ws.Range("P7").Value = 0.5989702197413518
$ws.Range("Q7").Value = 47.24733215916322
$ws.Range("R7").Value = 425.225989432469
$ws.Range("S7").Value = 0.03155636158701893
$ws.Range("T7").Value = 0.03155636158701893

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Mag"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 85.41274733333334
$ws.Range("H8").Value = 256.238242
$ws.Range("I8").Value = 0.04456888863285297
$ws.Range("J8").Value = 0.04456888863285297
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3133113333333333
$ws.Range("N8").Value = 0.939934
$ws.Range("O8").Value = 0.4010297802586483
$ws.Range("P8").Value = 0.4010297802586483
$ws.Range("Q8").Value = 26.76078175066978
$ws.Range("R8").Value = 240.847035756028
$ws.Range("S8").Value = 0.01787345161480519
$ws.Range("T8").Value = 0.01787345161480519

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Mag"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 85.41274733333334
$ws.Range("H9").Value = 256.238242
$ws.Range("I9").Value = 0.04456888863285297
$ws.Range("J9").Value = 0.04456888863285297
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4679556666666667
$ws.Range("N9").Value = 1.403867
$ws.Range("O9").Value = 0.5989702197413518
$ws.Range("P9").Value = 0.5989702197413518
$ws.Range("Q9").Value = 39.96937912020156
$ws.Range("R9").Value = 359.724412081814
$ws.Range("S9").Value = 0.02669543701804778
$ws.Range("T9").Value = 0.02669543701804778

